$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.058.75"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "1.581.50"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "'211.76"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  +6.97%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").Value = "'25.60"
$ws.Range("E8").Value = "  +8.94%  "
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "'0.0899"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.809.03"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "1.554.98"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "29.106.75"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'3.70"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.521"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "'62.41"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "'237.08"
$ws.Range("E18").Value = "  +4.53%  "
$ws.Range("D19").Value = "'7.41"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +5.52%  "
$ws.Range("D25").Value = "'153.11"
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("D26").Value = "'0.109"
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").Value = "'15.08"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "1.419.29"
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("D34").Value = "'3.04"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").Value = "'2.78"
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").Value = "'0.524"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'52.30"
$ws.Range("E43").Value = "  +24.26%  "
$ws.Range("D44").Value = "'0.786"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "'64.57"
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "1.719.23"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").Value = "'0.852"
$ws.Range("E49").Value = "  -6.04%  "
$ws.Range("D50").Value = "'85.57"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.0512"
$ws.Range("E51").Value = "  +0.73%  "
